{"js": "// The document body contains a single, single-column table whose rows hold\n// DaCapo/ZGC benchmark statistics (one value per row, as plain text).\n// This edit updates a handful of summary rows and replaces the raw,\n// tab-separated \"per-invocation\" rows (43-45) with their single\n// already-known summary value (the value that used to live in rows 0-2).\n\nconst table = context.document.body.tables.getFirstOrNullObject();\nawait context.sync();\n\ntable.load(\"values\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body but found none.\");\n}\n\n// rowIndex (0-based) -> new cell text\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"641\",\n  5: \"0.00096\",\n  6: \"0.00013\",\n  8: \"0.00045\",\n  9: \"0.00045\",\n  10: \"0.00045\",\n  11: \"0.03854\",\n  // These rows used to carry a single summary number, then a raw\n  // tab-separated dump of per-invocation numbers got appended into the\n  // same run; restore the single summary value.\n  43: \"99.95\",\n  44: \"0.04\",\n  45: \"72\",\n};\n\nfor (const [rowIndexStr, newValue] of Object.entries(updates)) {\n  const rowIndex = Number(rowIndexStr);\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = newValue;\n}\n\nawait context.sync();\n", "ps1": "# The document body contains a single, single-column table whose rows hold\n# DaCapo/ZGC benchmark statistics (one value per row, as plain text).\n# This edit updates a handful of summary rows and replaces the raw,\n# tab-separated \"per-invocation\" rows (table rows 44-46, 1-based) with\n# their single already-known summary value (the value that used to live\n# in rows 1-3).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row (1-based) -> new cell text\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"641\"\n    6  = \"0.00096\"\n    7  = \"0.00013\"\n    9  = \"0.00045\"\n    10 = \"0.00045\"\n    11 = \"0.00045\"\n    12 = \"0.03854\"\n    # These rows used to carry a single summary number, then a raw\n    # tab-separated dump of per-invocation numbers got appended into the\n    # same run; restore the single summary value.\n    44 = \"99.95\"\n    45 = \"0.04\"\n    46 = \"72\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
